$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Main "CERTIFICA" paragraph: update wording, month range and source
#    text ("banco de datos" -> "base de datos", months list, etc.)
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Que el comportamiento de la precipitación para el(los) mes(es) febrero, marzo, abril de 2021, de acuerdo con la información preliminar disponible en nuestro banco de datos de la estación meteorológica",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Que el comportamiento mensual de la precipitación y su respectivo índice I(%), durante el(los) mes(es) enero de 2021, 2023, de acuerdo con la información preliminar disponible en nuestra base de datos de la estación meteorológica",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2) Precipitation data table: add a 2021-01 row before 2021-02, fix the
#    2021-03 index value and append the rest of the monthly rows through
#    2022-02.
# ---------------------------------------------------------------------
$t = $d.Tables.Item(1)

# Insert the 2021-01 row right after the header row (before 2021-02).
$newRow = $t.Rows.Add($t.Rows.Item(2))
$newRow.Cells.Item(1).Range.Text = "2021-01"
$newRow.Cells.Item(2).Range.Text = "40.7"
$newRow.Cells.Item(3).Range.Text = "96.9"

# Fix the 2021-03 index value (134.88 -> 134.9). 2021-03 is now row 4.
$t.Cell(4, 3).Range.Text = "134.9"

# Append the remaining monthly rows (2021-04 .. 2022-02) at the end of
# the table, each with Fecha / Precipitación total mensual (mm) / Índice (%).
$rowsToAdd = @(
    @("2021-04", "296.6", "140.5"),
    @("2021-05", "377.4", "118.7"),
    @("2021-06", "298.0", "85.0"),
    @("2021-07", "320.8", "88.5"),
    @("2021-08", "317.8", "107.3"),
    @("2021-09", "185.7", "85.7"),
    @("2021-10", "158.7", "102.8"),
    @("2021-11", "63.3", "54.8"),
    @("2021-12", "12.100000000000001", "22.1"),
    @("2022-01", "14.399999999999999", "34.3"),
    @("2022-02", "47.9", "104.1")
)

foreach ($rowData in $rowsToAdd) {
    $r = $t.Rows.Add()
    $r.Cells.Item(1).Range.Text = $rowData[0]
    $r.Cells.Item(2).Range.Text = $rowData[1]
    $r.Cells.Item(3).Range.Text = $rowData[2]
}

# ---------------------------------------------------------------------
# 3) Remove the ">: Mayor que." and "<: Menor que." convention
#    paragraphs (no longer applicable to this template).
#    NOTE: use $d.Content.Paragraphs (not $d.Paragraphs) — after a table
#    row mutation the plain $d.Paragraphs collection can return stale
#    items, while re-deriving it from $d.Content stays accurate.
# ---------------------------------------------------------------------
$found = $true
while ($found) {
    $found = $false
    $paras = $d.Content.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $ptext = $paras.Item($i).Range.Text
        if ($ptext.StartsWith(">:") -or $ptext.StartsWith("<:")) {
            $paras.Item($i).Range.Delete()
            $found = $true
            break
        }
    }
}

# ---------------------------------------------------------------------
# 4) Update the certification date field text.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "28 de abril de 2024", $true, $false, $false, $false, $false, $true, 1, $false,
    "30 de mayo de 2024", 2) | Out-Null

# ---------------------------------------------------------------------
# 5) Normalize the signature line run/proofing-error boundaries around
#    "Teniente Coronel JORGE GIOVANNI JIMÉNEZ SÁNCHEZ" (merges the
#    split runs and clears the stale grammar-check markers).
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Teniente Coronel JORGE GIOVANNI JIMÉNEZ SÁNCHEZ ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Teniente Coronel JORGE GIOVANNI JIMÉNEZ SÁNCHEZ ",
    2) | Out-Null
